# handlocations/importfile_lokaties_PT.xlsx
# Add a name/label column in A, shift the existing measurement columns one
# to the right (B:O) and append 4 new rows of data (+ new row labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wipe the existing data + any column-width formatting on A:N so we start
# from a clean sheet (this also drops the custom widths that used to live
# on columns B:E, matching the target file which has no <cols> override).
$ws.Columns("A:N").Delete()

# Column A: person/location labels (now backed by shared strings).
$labels = @("daniel_l", "daniel_r", "paul_l", "paul_r", "rienco_l", "rienco_r", "thijs_l", "thijs_r")
for ($r = 0; $r -lt $labels.Length; $r++) {
    $ws.Cells.Item($r + 1, 1).Value = $labels[$r]
}

# Columns B:O: the 14 measurement values per row (rows 1-4 are the
# pre-existing numbers, rows 5-8 are newly added).
$row1 = @(0.51, 0.52, 0.4, 0.64, 0.56, 0.6, 0.85, 0.56, 0.58, 0.52, 0.72, 0.77, 1.28, 1.14)
$row2 = @(0.51, 0.53, 0.72, 0.24, 0.44, 0.69, 0.76, 0.44, 0.31, 0.44, 0.85, 1.1, 0.93, 1.25)
$row3 = @(0.37, 0.3, 0.36, 0.27, 0.29, 0.34, 0.35, 0.23, 0.16, 0.32, 0.37, 0.67, 0.38, 0.95)
$row4 = @(0.34, 0.15, 0.34, 0.26, 0.35, 0.42, 0.37, 0.26, 0.32, 0.38, 0.35, 0.42, 0.3, 0.4)
$row5 = @(0.55, 0.26, 0.38, 0.38, 0.4, 0.6, 0.48, 0.52, 0.42, 0.57, 0.59, 0.95, 0.98, 0.95)
$row6 = @(0.39, 0.29, 0.36, 0.52, 0.33, 0.73, 0.53, 0.45, 0.23, 0.52, 0.28, 0.98, 1.05, 0.86)
$row7 = @(0.28, 0.19, 0.39, 0.3, 0.24, 0.42, 0.2, 0.18, 0.36, 0.3, 0.3, 0.39, 0.4, 0.41)
$row8 = @(0.36, 0.24, 0.56, 0.48, 0.19, 0.32, 0.43, 0.22, 0.5, 0.14, 0.26, 0.3, 0.45, 0.43)
$data = @($row1, $row2, $row3, $row4, $row5, $row6, $row7, $row8)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowvals = $data[$r]
    for ($c = 0; $c -lt $rowvals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 2).Value = $rowvals[$c]
    }
}

# Selection moved to M17 and the workbook window was resized/repositioned.
$ws.Range("M17").Select()
$excel.ActiveWindow.WindowState = -4143
$excel.Left = 0
$excel.Top = 45
$excel.Width = 23955
$excel.Height = 10545
